# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column labels in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold / centered /
# bordered), by copying the existing header cell's format onto the new ones.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-53: team record values, the same for every player on the roster.
$firstRow = 2
$lastRow = 53
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 75   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 87   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
